$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.053.72'
$ws.Range('E2').Value = '  -0.95%  '
$ws.Range('D3').Value = '2.218.90'
$ws.Range('E3').Value = '  -1.68%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.60'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.627'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.89%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '73.29'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -4.28%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.606'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -3.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.44'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0957'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.00'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -4.70%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.103'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.29%  '
$ws.Range('D14').Value = '2.550.66'
$ws.Range('E14').Value = '  -1.72%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.26'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.836'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.79%  '
$ws.Range('D17').Value = '2.243.18'
$ws.Range('E17').Value = '  -0.51%  '
$ws.Range('D18').Value = '41.860.64'
$ws.Range('E18').Value = '  -1.09%  '
$ws.Range('E19').Value = '  +4.92%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.77'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.42%  '
$ws.Range('E21').Value = '  -0.97%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.84'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +17.51%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '229.81'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.00%  '
$ws.Range('E24').Value = '  -7.60%  '
$ws.Range('E25').Value = '  +0.70%  '
$ws.Range('E26').Value = '  +0.17%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.73'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +3.18%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.27'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.08%  '
$ws.Range('E29').Value = '  -1.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '167.69'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.81%  '
$ws.Range('E31').Value = '  -1.33%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.65'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +5.59%  '
$ws.Range('E33').Value = '  -4.27%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '30.11'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -3.23%  '
$ws.Range('E35').Value = '  -0.66%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.109'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -10.21%  '
$ws.Range('E37').Value = '  -6.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0301'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -5.29%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '13.86'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '65.13'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +2.53%  '
$ws.Range('E41').Value = '  -3.60%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.66'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -3.09%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.198'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -3.55%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.74'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '104.92'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -3.91%  '
$ws.Range('E46').Value = '  -2.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.37'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.68%  '
$ws.Range('E48').Value = '  -1.28%  '
$ws.Range('E49').Value = '  -2.41%  '
$ws.Range('E50').Value = '  +0.13%  '
$ws.Range('D51').Value = '2.424.70'
$ws.Range('E51').Value = '  -1.82%  '
